$p = $ppt.ActivePresentation

# --- Slide 23: "Dig That Lick Similarity Search" -----------------------
# Fix stray ", " -> "." typo in the first bullet.
$s23 = $p.Slides.Item(23)
$shp23 = $s23.Shapes.Item(2)
$tr23 = $shp23.TextFrame.TextRange
$run23 = $tr23.Characters(1, 50)
$run23.Text = "Allows similarity search for patterns in the WJD."

# The text edit above nudges this autofit text box's computed height;
# restore it to its originally authored size (unchanged by this commit).
$shp23.Height = 378.4046630859375

# --- Slide 25: "Feature History Explorer" -------------------------------
# Fix "What das " -> "What does " typo.
$s25 = $p.Slides.Item(25)
$shp25 = $s25.Shapes.Item(2)
$tr25 = $shp25.TextFrame.TextRange
$run25 = $tr25.Characters(39, 10)
$run25.Text = "What does "

# Grow the text box to its new authored height.
$shp25.Height = 394.6971740722656

# --- Slide 27: "Dig That Lick Pattern Search" ---------------------------
# Remove stray " Lick" from the question text.
$s27 = $p.Slides.Item(27)
$shp27 = $s27.Shapes.Item(2)
$tr27 = $shp27.TextFrame.TextRange
$run27 = $tr27.Characters(293, 180)
$run27.Text = "Repeat one of the previous searches with tone context of 2 or more tones before and after. How does it change the pattern impression? What are the most common pre/successions?"

# Shrink the text box to its new authored height.
$shp27.Height = 362.70782470703125
